$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "43.188.09"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +2.19%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.318.51"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +2.07%  "

$ws.Range("E4").Value = "  +0.03%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "303.64"
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "100.94"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +5.76%  "

$ws.Range("E7").Value = "  +2.73%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +3.59%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "34.78"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +4.09%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0798"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.12%  "

$ws.Range("E12").Value = "  +4.20%  "

$ws.Range("E13").Value = "  +15.64%  "

$ws.Range("E14").Value = "  +3.15%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "2.694.91"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.60%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "2.349.63"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +3.37%  "

$ws.Range("E17").Value = "  +5.34%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "43.145.90"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +2.25%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "12.62"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +8.46%  "

$ws.Range("E20").Value = "  +2.88%  "

$ws.Range("E21").Value = "  +1.90%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "67.92"
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "238.02"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.97%  "

$ws.Range("E24").Value = "  +12.60%  "

$ws.Range("E25").Value = "  +0.80%  "

$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("E27").Value = "  +4.29%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "168.40"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.48%  "

$ws.Range("E29").Value = "  -8.83%  "

$ws.Range("E30").Value = "  +0.57%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "9.23"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.31%  "

$ws.Range("E32").Value = "  +0.09%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.05"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.78%  "

$ws.Range("E34").Value = "  +4.01%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "17.25"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +5.16%  "

$ws.Range("E36").Value = "  +4.17%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.0696"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.71%  "

$ws.Range("E38").Value = "  +3.91%  "

$ws.Range("E39").Value = "  +4.48%  "

$ws.Range("E40").Value = "  +1.47%  "

$ws.Range("E41").Value = "  +0.69%  "

$ws.Range("E42").Value = "  -1.78%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.007.75"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.12%  "

$ws.Range("E44").Value = "  +3.18%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "10.15"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +5.71%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "17.72"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.87%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.87"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.80%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "56.03"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +7.20%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.537.47"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.58%  "

$ws.Range("E50").Value = "  +5.07%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "4.60"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.65%  "
